$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 0.8584576666666667
$ws.Range("N2").Value = 2.575373
$ws.Range("O2").Value = 0.02952026538348031
$ws.Range("P2").Value = 0.03028938521394646
$ws.Range("Q2").Value = 0.3334261023435556
$ws.Range("R2").Value = 3.000834921092
$ws.Range("S2").Value = 0.02952026538348031
$ws.Range("T2").Value = 0.03028938521394646

# Row 3
$ws.Range("N3").Value = 75.717583
$ws.Range("O3").Value = 0.8679143348771993
$ws.Range("P3").Value = 0.8905269407406087
$ws.Range("Q3").Value = 9.802936731325779
$ws.Range("R3").Value = 88.22643058193202
$ws.Range("S3").Value = 0.8679143348771993
$ws.Range("T3").Value = 0.8905269407406087

# Row 4
$ws.Range("M4").Value = 0.3580240000000001
$ws.Range("N4").Value = 1.074072
$ws.Range("O4").Value = 0.0123115721415754
$ws.Range("P4").Value = 0.01263233735676886
$ws.Range("Q4").Value = 0.1390569989653334
$ws.Range("R4").Value = 1.251512990688
$ws.Range("S4").Value = 0.0123115721415754
$ws.Range("T4").Value = 0.01263233735676886

# Row 5
$ws.Range("M5").Value = 2.2152535
$ws.Range("N5").Value = 4.430507
$ws.Range("O5").Value = 0.07617716487477769
$ws.Range("P5").Value = 0.05210792115009603
$ws.Range("Q5").Value = 0.8604074130713335
$ws.Range("R5").Value = 5.162444478428001
$ws.Range("S5").Value = 0.07617716487477769
$ws.Range("T5").Value = 0.05210792115009603

# Row 6
$ws.Range("M6").Value = 0.4093533333333334
$ws.Range("N6").Value = 1.22806
$ws.Range("O6").Value = 0.01407666272296744
$ws.Range("P6").Value = 0.01444341553857988
$ws.Range("Q6").Value = 0.1589933804711112
$ws.Range("R6").Value = 1.43094042424
$ws.Range("S6").Value = 0.01407666272296744
$ws.Range("T6").Value = 0.01444341553857988
